$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-11 Friday" "2024-10-12 Saturday"

Replace-Text "999×3=2997" "549×3=1647"
Replace-Text "447×8=3576" "849×9=7641"
Replace-Text "509×8=4072" "285×9=2565"
Replace-Text "673×6=4038" "456×8=3648"
Replace-Text "410×7=2870" "561×4=2244"

Replace-Text "821×6=4926" "670×5=3350"
Replace-Text "219×5=1095" "800×8=6400"
Replace-Text "391×4=1564" "283×6=1698"
Replace-Text "585×2=1170" "945×3=2835"
Replace-Text "441×9=3969" "880×6=5280"

Replace-Text "507×5=2535" "551×9=4959"
Replace-Text "563×6=3378" "221×6=1326"

Replace-Text "949×6=5694" "152×8=1216"
Replace-Text "948×6=5688" "655×7=4585"

Replace-Text "904×2=1808" "106×3=318"
Replace-Text "763×6=4578" "509×7=3563"
Replace-Text "953×7=6671" "551×2=1102"
Replace-Text "866×2=1732" "527×9=4743"
Replace-Text "332×9=2988" "361×6=2166"

Replace-Text "976×7=6832" "584×2=1168"
Replace-Text "345×4=1380" "396×7=2772"
Replace-Text "898×4=3592" "232×2=464"
Replace-Text "270×5=1350" "884×2=1768"
Replace-Text "174×6=1044" "526×6=3156"
